$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last two existing rows (24:25) into the two
# new rows (26:27) so the new cells reuse the same cell styles / shared
# strings as the rest of the table.
$ws.Range("A24:F25").Copy()
$ws.Range("A26:F27").PasteSpecial()

# Row 26 - 2025-10-13 (serial 45943), 四方坪站
$ws.Range("A26").Value = 45943
$ws.Range("C26").Value = 8434.09
$ws.Range("D26").Value = 6967.91
$ws.Range("E26").Value = 2913.27
$ws.Range("F26").Value = 369

# Row 27 - 2025-10-13 (serial 45943), 高岭站
$ws.Range("A27").Value = 45943
$ws.Range("C27").Value = 4333.8900000000003
$ws.Range("D27").Value = 3439.09
$ws.Range("E27").Value = 1107.18
$ws.Range("F27").Value = 169

$ws.Range("I27").Select()
